# Adding the Header to our Project
# 1) Remove the stray trailing <a:endParaRPr/> on the "Understanding Box Sizing" title (slide 40)
# 2) Insert a new slide ("Adding the Header to our Project") right before the final "Thanks!!!!" slide

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Step 1: slide 40 ("Understanding Box Sizing") - drop the empty endParaRPr
# that trails the title run so the title paragraph ends right after the run.
# ---------------------------------------------------------------------------
$boxSizingSlide = $p.Slides.Item(40)
$titleShape = $boxSizingSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Understanding Box Sizing"

# ---------------------------------------------------------------------------
# Step 2: new slide - duplicate slide 40 (same "Title and Content" layout)
# so the generated markup matches a normal, freshly authored slide (clean
# <a:lstStyle/>, no leftover autofit noise), then overwrite its content.
# ---------------------------------------------------------------------------
$dup = $boxSizingSlide.Duplicate()
$newSlide = $dup.Item(1)

# Title placeholder
$titleNew = $newSlide.Shapes.Item(1)
$titleNew.TextFrame.TextRange.Text = "Adding the Header to our Project"
$titleNew.TextFrame.TextRange.LanguageID = "en-IN"
$titleNew.Left = 53.33340582677165
$titleNew.Top = 15.27277590551181
$titleNew.Width = 676.9030121259842
$titleNew.Height = 51.00002

# Content placeholder
$bodyNew = $newSlide.Shapes.Item(2)
$bodyText = "Now lets remove the border too from the section.`rAlso I have added some code to our index.html to add an ugly looking navigation bar to the code.`rNow lets try to make the header look a bit better so first since we might use the header tag again in our code so lets add a class main-header to our header tag and use the class selector to add css to it`rSo now we want to have a navigation bar that spans the whole screen width ,has a green background color, and also it should have some padding so that the content of the header doesn’t sit on the edges`rSo now add width:100% , background:#2ddf5c,padding:8px 16px to achieve the full width,background color and padding.`r "
$bodyNew.TextFrame.TextRange.Text = $bodyText
$bodyNew.TextFrame.TextRange.LanguageID = "en-IN"
$bodyNew.Left = 53.33340582677165
$bodyNew.Top = 77.72734283464568
$bodyNew.Width = 864.6666341732283
$bodyNew.Height = 397.97057118110234

$paras = $bodyNew.TextFrame.TextRange.Paragraphs()

# Paragraph 1 & 2 are plain en-GB sentences
$paras.Item(1).LanguageID = "en-GB"
$paras.Item(2).LanguageID = "en-GB"

# Paragraph 3 is en-GB too
$paras.Item(3).LanguageID = "en-GB"

# Paragraphs 4-5 use en-IN (already set on whole body above)

# Last (6th) paragraph: single space run, no-bullet paragraph formatting
$lastPara = $paras.Item(6)
$lastPara.ParagraphFormat.Bullet.Visible = $false
$lastPara.ParagraphFormat.LeftIndent = 0
$lastPara.ParagraphFormat.FirstLineIndent = 0
